$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "classical-best-embed vs. classical-best-tfidf"
$ws.Range("C2").Value = 0.058
$ws.Range("E2").Value = 0.03
$ws.Range("F2").Value = 0.019
$ws.Range("H2").Value = 0.041
$ws.Range("I2").Value = 0.032
$ws.Range("J2").Value = 0.041

# Row 3
$ws.Range("C3").Value = 0.068
$ws.Range("D3").Value = 0.079
$ws.Range("E3").Value = 0.105
$ws.Range("F3").Value = 0.067
$ws.Range("G3").Value = 0.13
$ws.Range("H3").Value = 0.099
$ws.Range("I3").Value = 0.08

# Row 4
$ws.Range("A4").Value = "BERT-base vs. classical-best-embed"
$ws.Range("C4").Value = 0.01
$ws.Range("D4").Value = 0.056
$ws.Range("E4").Value = 0.075
$ws.Range("F4").Value = 0.048
$ws.Range("G4").Value = 0.053
$ws.Range("H4").Value = 0.058
$ws.Range("I4").Value = 0.047
$ws.Range("J4").Value = 0.05

# Row 5
$ws.Range("B5").Value = 0.338
$ws.Range("C5").Value = 0.158
$ws.Range("D5").Value = 0.076
$ws.Range("E5").Value = 0.094
$ws.Range("G5").Value = 0.045
$ws.Range("H5").Value = 0.016
$ws.Range("I5").Value = 0.091
$ws.Range("J5").Value = 0.071

# Row 6
$ws.Range("A6").Value = "BERT-base-nli vs. classical-best-embed"
$ws.Range("B6").Value = 0.338
$ws.Range("C6").Value = 0.1
$ws.Range("D6").Value = 0.053
$ws.Range("E6").Value = 0.064
$ws.Range("F6").Value = 0.018
$ws.Range("G6").Value = -0.032
$ws.Range("H6").Value = -0.025
$ws.Range("I6").Value = 0.059
$ws.Range("J6").Value = 0.03

# Row 7
$ws.Range("B7").Value = 0.338
$ws.Range("C7").Value = 0.09
$ws.Range("D7").Value = -0.003
$ws.Range("E7").Value = -0.011
$ws.Range("F7").Value = -0.03
$ws.Range("G7").Value = -0.085
$ws.Range("H7").Value = -0.083
$ws.Range("I7").Value = 0.011
$ws.Range("J7").Value = -0.02
